# Update the "userCount" column (F) values for the players table.
# Mirrors the diff: several F-column (userCount) values increased slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 66
    3  = 132
    4  = 92
    5  = 72
    6  = 468
    7  = 320
    8  = 79
    9  = 266
    10 = 96
    11 = 330
    12 = 49
    14 = 585
    15 = 103
    16 = 84
    17 = 215
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
